$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.931299999999998
$ws.Range("D21").Value = -7.536200000000001
$ws.Range("D23").Value = -6.287799999999993
$ws.Range("D25").Value = -8.465699999999998
